$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(42603.657893518517, "Named", 3312, 408, 16, 3, 8, 27, 72, 0, 1, 0, 100),
    @(42603.690972222219, "Named", 3121, 408, 16, 3, 8, 27, 72, 0, 1, 0, 100),
    @(42603.691527777781, "Named", 3164, 408, 16, 3, 8, 27, 72, 0, 1, 0, 100),
    @(42603.692523148151, "Named", 3118, 408, 16, 3, 8, 27, 72, 0, 1, 0, 100),
    @(42603.692824074074, "Named", 3034, 408, 16, 3, 8, 27, 72, 0, 1, 0, 100),
    @(42603.693449074075, "Named", 3182, 408, 16, 3, 8, 27, 72, 0, 1, 0, 100)
)

$startRow = 33
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($col = 0; $col -lt $rowData.Count; $col++) {
        $ws.Cells.Item($row, $col + 1).Value = $rowData[$col]
    }
}
